$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 664
$ws.Range("G6").Value = 88
$ws.Range("G7").Value = 40
$ws.Range("F8").Value = 2713
$ws.Range("F10").Value = 6365
$ws.Range("F11").Value = 2402
$ws.Range("F15").Value = 2570
$ws.Range("F16").Value = 31
$ws.Range("F18").Value = 6811
$ws.Range("F19").Value = 250
$ws.Range("F21").Value = 185
$ws.Range("F24").Value = 7689
$ws.Range("F36").Value = 46
$ws.Range("F37").Value = 65
$ws.Range("F38").Value = 2573
$ws.Range("F42").Value = 1147
$ws.Range("F44").Value = 596
$ws.Range("F45").Value = 3607
$ws.Range("F46").Value = 134
$ws.Range("F48").Value = 95
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 25
$ws.Range("F6").Value = 3
$ws.Range("F15").Value = 163
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 664
$ws.Range("G4").Value = 88
$ws.Range("G5").Value = 40
$ws.Range("F7").Value = 2713
$ws.Range("F8").Value = 25
$ws.Range("F10").Value = 6365
$ws.Range("F11").Value = 2402
$ws.Range("F15").Value = 2570
$ws.Range("F16").Value = 31
$ws.Range("F20").Value = 6811
$ws.Range("F21").Value = 250
$ws.Range("F23").Value = 185
$ws.Range("F25").Value = 7689
$ws.Range("F35").Value = 46
$ws.Range("F37").Value = 65
$ws.Range("F38").Value = 2573
$ws.Range("F41").Value = 1147
$ws.Range("F43").Value = 596
$ws.Range("F44").Value = 163
$ws.Range("F45").Value = 3607
$ws.Range("F46").Value = 134
$ws.Range("F49").Value = 95
